$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.877.35'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '1.743.00'
$ws.Range("E3").Value = '  -0.37%  '

$ws.Range("D4").Formula = "'1.002"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Formula = "'230.82"
$ws.Range("E5").Value = '  -2.19%  '

$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").Formula = "'0.5162"
$ws.Range("E7").Value = '  +1.52%  '

$ws.Range("D8").Formula = "'0.2793"
$ws.Range("E8").Value = '  +4.83%  '

$ws.Range("D9").Formula = "'39.46"
$ws.Range("E9").Value = '  -3.22%  '

$ws.Range("D10").Formula = "'0.06097"
$ws.Range("E10").Value = '  -1.52%  '

$ws.Range("D11").Value = '1.752.55'
$ws.Range("E11").Value = '  -0.06%  '

$ws.Range("D12").Formula = "'0.07040"
$ws.Range("E12").Value = '  +1.60%  '

$ws.Range("D13").Formula = "'15.23"
$ws.Range("E13").Value = '  -1.04%  '

$ws.Range("D14").Formula = "'0.6388"
$ws.Range("E14").Value = '  +1.95%  '

$ws.Range("D15").Formula = "'4.503"
$ws.Range("E15").Value = '  +0.67%  '

$ws.Range("D16").Formula = "'77.05"
$ws.Range("E16").Value = '  -0.73%  '

$ws.Range("D17").Formula = "'1.002"
$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("E18").Value = '  -0.01%  '

$ws.Range("D19").Value = '25.869.96'
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Formula = "'11.43"
$ws.Range("E20").Value = '  -1.67%  '

$ws.Range("D21").Formula = "'0.000006565"
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").Value = '1.973.88'
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").Formula = "'4.129"
$ws.Range("E23").Value = '  +1.86%  '

$ws.Range("D24").Formula = "'8.624"
$ws.Range("E24").Value = '  +4.48%  '

$ws.Range("D25").Formula = "'5.132"
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Formula = "'139.69"
$ws.Range("E26").Value = '  +2.22%  '

$ws.Range("D27").Formula = "'1.514"
$ws.Range("E27").Value = '  +4.01%  '

$ws.Range("D28").Formula = "'15.05"
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").Formula = "'1.805"
$ws.Range("E29").Value = '  +3.66%  '

$ws.Range("D30").Formula = "'102.26"
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("D31").Formula = "'0.08227"
$ws.Range("E31").Value = '  +0.44%  '

$ws.Range("D32").Formula = "'3.654"
$ws.Range("E32").Value = '  -1.11%  '

$ws.Range("D33").Formula = "'3.419"
$ws.Range("E33").Value = '  +0.57%  '

$ws.Range("D34").Formula = "'0.04484"
$ws.Range("E34").Value = '  +1.53%  '

$ws.Range("D35").Formula = "'2.617"
$ws.Range("E35").Value = '  -1.29%  '

$ws.Range("D36").Formula = "'0.9785"
$ws.Range("E36").Value = '  -1.71%  '

$ws.Range("D37").Formula = "'0.6128"
$ws.Range("E37").Value = '  +2.39%  '

$ws.Range("D38").Formula = "'2.648"
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("D39").Formula = "'0.01586"
$ws.Range("E39").Value = '  +1.51%  '

$ws.Range("D40").Formula = "'1.919"
$ws.Range("E40").Value = '  -1.14%  '

$ws.Range("D41").Formula = "'1.002"
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").Formula = "'100.62"
$ws.Range("E42").Value = '  -0.71%  '

$ws.Range("D43").Formula = "'0.3826"
$ws.Range("E43").Value = '  +0.23%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Formula = "'0.7208"
$ws.Range("E44").Value = '  -3.89%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Formula = "'4.972"
$ws.Range("E45").Value = '  +1.70%  '

$ws.Range("D46").Formula = "'0.05422"
$ws.Range("E46").Value = '  -1.48%  '

$ws.Range("D47").Formula = "'6.258"
$ws.Range("E47").Value = '  +5.61%  '

$ws.Range("D48").Formula = "'0.1120"
$ws.Range("E48").Value = '  +2.15%  '

$ws.Range("D49").Formula = "'53.17"
$ws.Range("E49").Value = '  +0.93%  '

$ws.Range("D50").Formula = "'7.673"
$ws.Range("E50").Value = '  +3.16%  '

$ws.Range("D51").Formula = "'29.79"
$ws.Range("E51").Value = '  -0.82%  '
